$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.302.73"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "3.116.65"
$ws.Range("E3").Value = "  -5.34%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'213.09"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").Value = "'628.02"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.393"
$ws.Range("E7").Value = "  -5.55%  "
$ws.Range("D8").Value = "'0.724"
$ws.Range("E8").Value = "  +2.28%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.116.90"
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "'0.553"
$ws.Range("E11").Value = "  -6.84%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  -5.96%  "
$ws.Range("D14").Value = "89.000.96"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "'5.25"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "3.690.79"
$ws.Range("E16").Value = "  -5.11%  "
$ws.Range("D17").Value = "'31.97"
$ws.Range("E17").Value = "  -6.99%  "
$ws.Range("D18").Value = "3.142.70"
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "'3.32"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").Value = "'0.0000212"
$ws.Range("E20").Value = "  +17.00%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("E21").Value = "  -7.52%  "
$ws.Range("D22").Value = "'424.61"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Value = "'8.31"
$ws.Range("E23").Value = "  -7.52%  "
$ws.Range("D24").Value = "'4.87"
$ws.Range("E24").Value = "  -9.03%  "
$ws.Range("D25").Value = "'5.25"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("E26").Value = "  -6.71%  "
$ws.Range("D27").Value = "'79.08"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "3.351.80"
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'0.157"
$ws.Range("E31").Value = "  -8.91%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.20"
$ws.Range("E32").Value = "  -5.61%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "'3.88"
$ws.Range("E33").Value = "  +5.60%  "
$ws.Range("D34").Value = "'507.95"
$ws.Range("E34").Value = "  -8.98%  "
$ws.Range("D35").Value = "'6.79"
$ws.Range("E35").Value = "  -6.80%  "
$ws.Range("D36").Value = "'1.84"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("E37").Value = "  -7.00%  "
$ws.Range("D38").Value = "'21.83"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").Value = "'22.18"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -6.56%  "
$ws.Range("E43").Value = "  -7.01%  "
$ws.Range("D44").Value = "'0.363"
$ws.Range("E44").Value = "  -7.66%  "
$ws.Range("D45").Value = "'145.01"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D46").Value = "'43.58"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "'165.84"
$ws.Range("E47").Value = "  -10.01%  "
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").Value = "'0.720"
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").Value = "'24.35"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("E51").Value = "  -8.60%  "

# Normalize style so forced-text cells don't retain a quote-prefix style marker
$ws.Range("D2:D51").Style = "Normal"
